$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("death_statistics")

# --- "death_statistics": add a "gamma recovery" column next to the
# existing "gamma" (death) column, and add the new Ferguson et al. priors
# for both death and recovery (mean/sd), per commit message:
# "added ferguson priors for death and recovery, gamma mean and sd"

# 1) Insert a new column before the existing "gamma" column. This shifts
#    the old B column (gamma / mean=15,sd=7 / mean=19,sd=9) into column C
#    and keeps the existing hyperlinks on A2/A3 intact.
$ws.Columns("B:B").Insert()

# 2) Headers (rename the existing "gamma" column to "gamma death" first,
#    then add the brand-new "gamma recovery" column).
$ws.Range("C1").Value = "gamma death"
$ws.Range("B1").Value = "gamma recovery"

# 3) Row 2 (NCBI death source, unchanged) - just clear the stray blank
#    cell left behind by the column insert in B2.
$ws.Range("B2").Clear()

# 4) Row 3: replace the old Imperial College row with the new Ferguson
#    medRxiv source, and fill in the death (C3) / recovery (B3) stats.
$ws.Range("A3").Value = "https://www.medrxiv.org/content/10.1101/2020.03.09.20033357v1.full.pdf"
$ws.Range("C3").Value = "mean = 18.8, sd = 8.5"
$ws.Range("B3").Value = "mean = 24.7, sd = 8.7"
# The column insert left B3 carrying a stray "Hyperlink" style copied
# from A3; reset it back to the sheet's plain (shrink-to-fit) style.
$ws.Range("B3").Style = "Normal"
$ws.Range("B3").ShrinkToFit = $true

# 5) Row 4: brand-new source row (plain text URL, no hyperlink styling).
$ws.Range("A4").Value = "https://www.medrxiv.org/content/10.1101/2020.04.01.20050138v1.full.pdf"

# --- Hyperlinks -------------------------------------------------------
# This runtime's Hyperlinks collection can only be rebuilt wholesale, so
# drop all existing links on the sheet and recreate rId1 (unchanged NCBI
# link) and rId2 (now pointing at the new medRxiv source) in order.
$ws.Range("A3").Hyperlinks.Delete()

$null = $ws.Hyperlinks.Add($ws.Range("A2"), "https://www.ncbi.nlm.nih.gov/pmc/articles/PMC7074197/")
$null = $ws.Hyperlinks.Add($ws.Range("A3"), "https://www.medrxiv.org/content/10.1101/2020.03.09.20033357v1.full.pdf")

# Restore the cell text-formatting that existed before the relink (Add()
# resets it to a freshly minted "Hyperlink" style): A2 keeps its original
# shrink-to-fit hyperlink look, A3 gets the plain hyperlink look.
$ws.Range("A2").ShrinkToFit = $true
$ws.Range("A3").Style = "Hyperlink"

# --- Selection ----------------------------------------------------------
$ws.Range("A5").Select()
